$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Standorte")

$ws.Range("C2").Value = 200
$ws.Range("C3").Value = 150
$ws.Range("C4").Value = 120
$ws.Range("C5").Value = 180
